# Splits the single run of text in a paragraph into one run per word and
# one run per inter-word space, matching the target OOXML diff, which
# turns a single <w:r><w:t>...</w:t></w:r> into an alternating sequence of
# <w:r> elements: word, space, word, space, ... word.

function Split-ParagraphIntoWordRuns($Document, $Paragraph) {
    $pRange = $Paragraph.Range
    $fullText = $pRange.Text
    # Drop the trailing paragraph mark (and any sectPr/cell markers) so we
    # only touch the paragraph's actual textual content.
    $mark = [char]13
    $cut = $fullText.IndexOf($mark)
    if ($cut -ge 0) {
        $text = $fullText.Substring(0, $cut)
    } else {
        $text = $fullText
    }

    # Build the target Range that spans just the text (not the pilcrow).
    $startPos = $pRange.Start
    $endPos = $startPos + $text.Length
    $targetRange = $Document.Range($startPos, $endPos)

    # Split the text on single spaces, re-inserting a standalone space run
    # between each pair of words.
    $words = $text.Split(" ")

    $sb = New-Object System.Text.StringBuilder
    [void]$sb.Append('<?xml version="1.0" encoding="UTF-8" standalone="yes"?>')
    [void]$sb.Append('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">')
    [void]$sb.Append('<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">')
    [void]$sb.Append('<pkg:xmlData>')
    [void]$sb.Append('<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">')
    [void]$sb.Append('<w:body><w:p>')

    for ($i = 0; $i -lt $words.Length; $i++) {
        if ($i -gt 0) {
            [void]$sb.Append('<w:r><w:t xml:space="preserve"> </w:t></w:r>')
        }
        $word = $words[$i]
        $escaped = $word.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        [void]$sb.Append('<w:r><w:t xml:space="preserve">')
        [void]$sb.Append($escaped)
        [void]$sb.Append('</w:t></w:r>')
    }

    [void]$sb.Append('</w:p></w:body></w:document>')
    [void]$sb.Append('</pkg:xmlData></pkg:part></pkg:package>')

    $targetRange.InsertXML($sb.ToString())
}

$d = $word.ActiveDocument

# Locate the Title / Author / Abstract paragraphs by their paragraph style
# (rather than a hard-coded index) so the edit is resilient to unrelated
# structural changes elsewhere in the document.
$targetStyles = @("Title", "Author", "Abstract")
$found = @{}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $styleName = $para.Style.NameLocal
    foreach ($style in $targetStyles) {
        if (-not $found.ContainsKey($style) -and $styleName -eq $style) {
            $found[$style] = $para
        }
    }
}

foreach ($style in $targetStyles) {
    if ($found.ContainsKey($style)) {
        Split-ParagraphIntoWordRuns $d $found[$style]
    }
}
